$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pBdr = "<w:pBdr><w:top w:val=`"single`" w:sz=`"4`" w:space=`"1`" w:color=`"auto`"/><w:left w:val=`"single`" w:sz=`"4`" w:space=`"4`" w:color=`"auto`"/><w:bottom w:val=`"single`" w:sz=`"4`" w:space=`"1`" w:color=`"auto`"/><w:right w:val=`"single`" w:sz=`"4`" w:space=`"4`" w:color=`"auto`"/></w:pBdr>"

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so that paragraph indices for
# content above the current edit point remain stable across edits.
# ---------------------------------------------------------------------------

# 6) Last paragraph ("-Proporcionar al cliente...") -> keep its content and
#    append a new paragraph right after it (replace the whole range, since
#    InsertXML on a collapsed range at the very end of the body clobbers the
#    final paragraph instead of appending after it).
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$xmlLast = "<w:p $wns><w:r><w:t>-Proporcionar al cliente un número de reserva, y que aparezca en pantalla un mensaje confirmando su reserva, o mejor aún, que esto se le envíe al email proporcionado.</w:t></w:r></w:p>" +
  "<w:p $wns><w:r><w:t>-Separar el código de routes/web.php mediante un include?</w:t></w:r></w:p>"
$rLast.InsertXML($xmlLast)

# 5) "¿Qué me falta por hacer?" paragraph -> drop the <w:lastRenderedPageBreak/>
#    (it moves to the newly-inserted "Archivos de configuración..." paragraph).
$pQue = $d.Paragraphs.Item(13)
$rQue = $pQue.Range
$xmlQue = "<w:p $wns><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val=`"yellow`"/></w:rPr><w:t>¿Qué me falta por hacer?</w:t></w:r></w:p>"
$rQue.InsertXML($xmlQue)

# 4) Empty paragraph right after the image (paragraph 12) -> stays, but gains
#    six new paragraphs right after it: four with content, one blank
#    paragraph carrying only border+bold pPr formatting, then one truly
#    blank <w:p/>. InsertXML replaces the whole range's contents, so the
#    original blank paragraph must be re-emitted explicitly.
$pBlank = $d.Paragraphs.Item(12)
$rBlank = $pBlank.Range
$xmlInsert =
  "<w:p $wns/>" +
  "<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t>Archivos de configuración para el envío automático del email:</w:t></w:r></w:p>" +
  "<w:p $wns><w:pPr>$pBdr<w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>app\Mail\ReservaConfirmada.php</w:t></w:r><w:r><w:t>: la clase del email</w:t></w:r></w:p>" +
  "<w:p $wns><w:pPr>$pBdr</w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>resources\views\emails/reserva_confirmada.blade.php</w:t></w:r><w:r><w:t>: la pantilla markdown del email</w:t></w:r></w:p>" +
  "<w:p $wns><w:pPr>$pBdr</w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-Routes/web.php</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space=`"preserve`"> Este interviene siempre en todos</w:t></w:r></w:p>" +
  "<w:p $wns><w:pPr>$pBdr<w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>" +
  "<w:p $wns/>"
$rBlank.InsertXML($xmlInsert)

# 3) Image paragraph (11) -> add <w:noProof/> ahead of the existing border in rPr.
$pImg = $d.Paragraphs.Item(11)
$pImg.Range.Font.NoProofing = -1

# 2) Merge the two runs "-resources/views/" + "inicio.blade.p" (paragraph 5) into one run.
$pInicio = $d.Paragraphs.Item(5)
$rInicio = $pInicio.Range
$xmlInicio = "<w:p $wns><w:pPr>$pBdr</w:pPr><w:r><w:t>-resources/views/inicio.blade.p</w:t></w:r><w:r><w:t>hp</w:t></w:r></w:p>"
$rInicio.InsertXML($xmlInicio)

# 1) "-Routes/web.php" paragraph (3) -> drop the <w:lang w:val="de-DE"/> formatting
#    from both the paragraph mark run properties and the text run.
$pRoutes = $d.Paragraphs.Item(3)
$rRoutes = $pRoutes.Range
$xmlRoutes = "<w:p $wns><w:pPr>$pBdr</w:pPr><w:r><w:t>-Routes/web.php</w:t></w:r></w:p>"
$rRoutes.InsertXML($xmlRoutes)

Write-Output "done"
